$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Indicator Data")

# Column B changes from a text "Survey" label (e.g. "2005-06 DHS") to a plain
# numeric year (the later year of the survey range), for every data row.
$years = @{
    3 = 2015; 4 = 2010; 5 = 2005; 6 = 2000; 7 = 1995; 8 = 1990; 9 = 1986;
    10 = 2006; 11 = 1999; 12 = 1993;
    13 = 2012; 14 = 2007; 15 = 2003; 16 = 1997; 17 = 1994; 18 = 1991; 19 = 1987;
    20 = 2015; 21 = 2014; 22 = 2009; 23 = 2003; 24 = 1998; 25 = 1993; 26 = 1989;
    27 = 2016; 28 = 2015; 29 = 2014; 30 = 2013; 31 = 2011; 32 = 2009; 33 = 2006; 34 = 2005; 35 = 1997; 36 = 1993; 37 = 1986
}

foreach ($row in $years.Keys) {
    $ws.Cells.Item($row, 2).Value = $years[$row]
}

# Restore the view/selection state recorded for this edit: frozen pane scrolled
# to B3, with B14 as the active selection in the bottom-right pane.
$ws.Range("B14").Select()
